$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is numeric-looking text (e.g. price strings like '596.83')
# must be explicitly formatted as Text first, otherwise Excel auto-converts them
# to numbers and loses exact formatting / introduces floating-point drift.
$textCells = @(
    'D5'
    'D6'
    'D13'
    'D16'
    'D20'
    'D25'
    'D28'
    'D31'
    'D34'
    'D38'
    'D40'
    'D41'
    'D44'
    'D46'
    'D48'
    'D49'
    'D50'
    'D51'
)
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$updates = @(
    @{ Cell = 'D2'; Value = '67.766.21' }
    @{ Cell = 'E2'; Value = '  +0.59%  ' }
    @{ Cell = 'D3'; Value = '3.802.21' }
    @{ Cell = 'E3'; Value = '  +0.59%  ' }
    @{ Cell = 'E4'; Value = '  +0.05%  ' }
    @{ Cell = 'D5'; Value = '596.83' }
    @{ Cell = 'E5'; Value = '  +0.60%  ' }
    @{ Cell = 'D6'; Value = '167.34' }
    @{ Cell = 'E6'; Value = '  +0.74%  ' }
    @{ Cell = 'E7'; Value = '  -0.11%  ' }
    @{ Cell = 'E9'; Value = '  +1.48%  ' }
    @{ Cell = 'E10'; Value = '  -1.20%  ' }
    @{ Cell = 'E11'; Value = '  +0.29%  ' }
    @{ Cell = 'E12'; Value = '  -0.45%  ' }
    @{ Cell = 'D13'; Value = '35.91' }
    @{ Cell = 'E13'; Value = '  +0.17%  ' }
    @{ Cell = 'D14'; Value = '4.443.44' }
    @{ Cell = 'D15'; Value = '3.801.79' }
    @{ Cell = 'E15'; Value = '  +0.52%  ' }
    @{ Cell = 'D16'; Value = '18.57' }
    @{ Cell = 'E16'; Value = '  +2.45%  ' }
    @{ Cell = 'D17'; Value = '67.810.30' }
    @{ Cell = 'E18'; Value = '  +1.59%  ' }
    @{ Cell = 'E19'; Value = '  +0.68%  ' }
    @{ Cell = 'D20'; Value = '460.87' }
    @{ Cell = 'E20'; Value = '  +0.74%  ' }
    @{ Cell = 'E21'; Value = '  -3.00%  ' }
    @{ Cell = 'E22'; Value = '  -0.05%  ' }
    @{ Cell = 'E23'; Value = '  +1.11%  ' }
    @{ Cell = 'E24'; Value = '  -0.08%  ' }
    @{ Cell = 'D25'; Value = '12.08' }
    @{ Cell = 'E25'; Value = '  +2.25%  ' }
    @{ Cell = 'E26'; Value = '  -1.23%  ' }
    @{ Cell = 'E27'; Value = '  +0.00%  ' }
    @{ Cell = 'D28'; Value = '10.00' }
    @{ Cell = 'E28'; Value = '  +0.57%  ' }
    @{ Cell = 'D29'; Value = '3.948.78' }
    @{ Cell = 'E29'; Value = '  +0.44%  ' }
    @{ Cell = 'E30'; Value = '  -0.35%  ' }
    @{ Cell = 'D31'; Value = '7.38' }
    @{ Cell = 'E31'; Value = '  +2.86%  ' }
    @{ Cell = 'E32'; Value = '  +2.03%  ' }
    @{ Cell = 'E33'; Value = '  -0.76%  ' }
    @{ Cell = 'D34'; Value = '0.999' }
    @{ Cell = 'E34'; Value = '  +0.01%  ' }
    @{ Cell = 'E35'; Value = '  -1.13%  ' }
    @{ Cell = 'D36'; Value = '3.742.67' }
    @{ Cell = 'E36'; Value = '  +0.23%  ' }
    @{ Cell = 'E37'; Value = '  +0.43%  ' }
    @{ Cell = 'D38'; Value = '3.35' }
    @{ Cell = 'E38'; Value = '  +1.94%  ' }
    @{ Cell = 'E39'; Value = '  +0.05%  ' }
    @{ Cell = 'D40'; Value = '0.999' }
    @{ Cell = 'E40'; Value = '  +0.85%  ' }
    @{ Cell = 'D41'; Value = '5.78' }
    @{ Cell = 'E41'; Value = '  +1.26%  ' }
    @{ Cell = 'E42'; Value = '  -0.09%  ' }
    @{ Cell = 'D44'; Value = '48.10' }
    @{ Cell = 'E44'; Value = '  +2.37%  ' }
    @{ Cell = 'E45'; Value = '  +1.53%  ' }
    @{ Cell = 'D46'; Value = '42.81' }
    @{ Cell = 'E46'; Value = '  -2.19%  ' }
    @{ Cell = 'E47'; Value = '  +0.07%  ' }
    @{ Cell = 'B48'; Value = 'Monero' }
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr' }
    @{ Cell = 'D48'; Value = '147.63' }
    @{ Cell = 'E48'; Value = '  -0.29%  ' }
    @{ Cell = 'B49'; Value = 'ONDO' }
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo' }
    @{ Cell = 'D49'; Value = '1.36' }
    @{ Cell = 'E49'; Value = '  +9.35%  ' }
    @{ Cell = 'B50'; Value = 'EnergySwap' }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens' }
    @{ Cell = 'D50'; Value = '27.22' }
    @{ Cell = 'E50'; Value = '  +7.51%  ' }
    @{ Cell = 'D51'; Value = '395.64' }
    @{ Cell = 'E51'; Value = '  +0.69%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
